$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.607.48'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').Value = '3.392.95'
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.79'
$ws.Range('E5').Value = '  +0.90%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.96'
$ws.Range('E6').Value = '  +0.44%  '

$ws.Range('E8').Value = '  -0.26%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.62'
$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('E10').Value = '  -0.50%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.385'
$ws.Range('E11').Value = '  -0.96%  '

$ws.Range('D12').Value = '3.975.87'
$ws.Range('E12').Value = '  +0.25%  '

$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.04'
$ws.Range('E14').Value = '  +0.95%  '

$ws.Range('D15').Value = '3.397.91'
$ws.Range('E15').Value = '  +0.25%  '

$ws.Range('E16').Value = '  -0.75%  '

$ws.Range('D17').Value = '61.651.28'
$ws.Range('E17').Value = '  +0.96%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.13'
$ws.Range('E18').Value = '  +0.73%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.64'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.17'
$ws.Range('E20').Value = '  +1.98%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.99'
$ws.Range('E21').Value = '  +1.24%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.45'
$ws.Range('E22').Value = '  -0.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.548'
$ws.Range('E23').Value = '  -0.77%  '

$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000114'
$ws.Range('E25').Value = '  -2.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.181'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.10%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.38'
$ws.Range('E28').Value = '  +1.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.98'
$ws.Range('E29').Value = '  -0.33%  '

$ws.Range('E30').Value = '  -0.38%  '

$ws.Range('E31').Value = '  +0.13%  '

$ws.Range('E32').Value = '  -0.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.33'

$ws.Range('E34').Value = '  -0.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '168.89'
$ws.Range('E35').Value = '  +0.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.11'
$ws.Range('E36').Value = '  +1.98%  '

$ws.Range('D37').Value = '3.425.68'
$ws.Range('E37').Value = '  +0.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.47'
$ws.Range('E38').Value = '  -0.09%  '

$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.43'
$ws.Range('E39').Value = '  +2.25%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0758'
$ws.Range('E40').Value = '  -1.29%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.781'
$ws.Range('E41').Value = '  +0.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.42'
$ws.Range('E42').Value = '  +0.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.66'

$ws.Range('E44').Value = '  +2.24%  '

$ws.Range('D45').Value = '2.487.72'
$ws.Range('E45').Value = '  +1.36%  '

$ws.Range('E46').Value = '  -1.33%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.64'
$ws.Range('E47').Value = '  -1.24%  '

$ws.Range('E48').Value = '  +0.01%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0264'
$ws.Range('E49').Value = '  -0.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.03'
$ws.Range('E50').Value = '  -6.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.205'
$ws.Range('E51').Value = '  -1.36%  '
